$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New averaged-intensity data for rows 10-19 (columns B:M), covering the three
# newly-added spiral sampling schemes plus the refreshed values for the
# schemes that follow them in the list.
$data = New-Object 'object[,]' 10,12
$data[0,0] = "Gaussian-Quadrature"
$data[0,1] = 1.632336121838838
$data[0,2] = 1.039049773750829
$data[0,3] = 0.8669497803913691
$data[0,4] = 1.632336121838838
$data[0,5] = 0.9407265842407487
$data[0,6] = 1.56056169696832
$data[0,7] = 0.9008325933958826
$data[0,8] = 1.039049773750829
$data[0,9] = 0.9529997770710992
$data[0,10] = 1.292667949454969
$data[0,11] = 1.156742758430998
$data[1,0] = "Spiral-90deg-10rot-5space"
$data[1,1] = 0.865565168152822
$data[1,2] = 1.779041957209686
$data[1,3] = 1.045593722660555
$data[1,4] = 0.865565168152822
$data[1,5] = 0.6339751667770079
$data[1,6] = 2.686634994849316
$data[1,7] = 0.825525129067891
$data[1,8] = 1.779041957209686
$data[1,9] = 1.41231783993512
$data[1,10] = 1.138941504043971
$data[1,11] = 1.306056023119546
$data[2,0] = "Spiral-90deg-15rot-5space"
$data[2,1] = 0.8671073559364372
$data[2,2] = 1.78440019139405
$data[2,3] = 1.043634011371434
$data[2,4] = 0.8671073559364372
$data[2,5] = 0.6362428823114351
$data[2,6] = 2.671690148825832
$data[2,7] = 0.8245356890149779
$data[2,8] = 1.78440019139405
$data[2,9] = 1.414017101382742
$data[2,10] = 1.14056222865959
$data[2,11] = 1.304601713142361
$data[3,0] = "Spiral-90deg-10rot-3space"
$data[3,1] = 0.8651370559443867
$data[3,2] = 1.780764868110878
$data[3,3] = 1.042980660867788
$data[3,4] = 0.8651370559443867
$data[3,5] = 0.6350291646224996
$data[3,6] = 2.689090987602157
$data[3,7] = 0.8233591244772326
$data[3,8] = 1.780764868110878
$data[3,9] = 1.411872764489333
$data[3,10] = 1.13850491021686
$data[3,11] = 1.306060310270824
$data[4,0] = "NoRotation-tilt60deg"
$data[4,1] = 0.4367920000000012
$data[4,2] = 1.273840000000003
$data[4,3] = 2.231227999999997
$data[4,4] = 0.4367920000000012
$data[4,5] = 0.4703840000000009
$data[4,6] = 1.422536000000001
$data[4,7] = 1.358055999999997
$data[4,8] = 1.273840000000003
$data[4,9] = 1.752534
$data[4,10] = 1.094663000000001
$data[4,11] = 1.198806
$data[5,0] = "Rotation-NoTilt"
$data[5,1] = 0.02
$data[5,2] = 0
$data[5,3] = 3.503762500000005
$data[5,4] = 0.02
$data[5,5] = 0.07000000000000001
$data[5,6] = 0.8798625000000005
$data[5,7] = 1.980837499999997
$data[5,8] = 0
$data[5,9] = 1.751881250000002
$data[5,10] = 0.8859406250000011
$data[5,11] = 1.07574375
$data[6,0] = "Rotation-60detTilt"
$data[6,1] = 0.4237672132608021
$data[6,2] = 0.3769460998144054
$data[6,3] = 2.380215189708796
$data[6,4] = 0.4237672132608021
$data[6,5] = 0.475242249113601
$data[6,6] = 0.9690109671424031
$data[6,7] = 1.531139587072001
$data[6,8] = 0.3769460998144054
$data[6,9] = 1.378580644761601
$data[6,10] = 0.9011739290112014
$data[6,11] = 1.026053551018668
$data[7,0] = "HexGrid-90degTilt5degRes"
$data[7,1] = 0.9854415717608653
$data[7,2] = 0.996290706389217
$data[7,3] = 0.9956980095968193
$data[7,4] = 0.9854415717608653
$data[7,5] = 0.9881389825467225
$data[7,6] = 0.9998535884593858
$data[7,7] = 0.9856501751847209
$data[7,8] = 0.996290706389217
$data[7,9] = 0.9959943579930182
$data[7,10] = 0.9907179648769417
$data[7,11] = 0.9918455056562885
$data[8,0] = "HexGrid-90degTilt22p5degRes"
$data[8,1] = 1.124001737165761
$data[8,2] = 1.084124609181023
$data[8,3] = 0.9738166222892486
$data[8,4] = 1.124001737165761
$data[8,5] = 0.9667899078408629
$data[8,6] = 0.8942063640296254
$data[8,7] = 0.9737782363777726
$data[8,8] = 1.084124609181023
$data[8,9] = 1.028970615735136
$data[8,10] = 1.076486176450449
$data[8,11] = 1.002786246147383
$data[9,0] = "HexGrid-60degTilt5degRes"
$data[9,1] = 0.9250863404092976
$data[9,2] = 0.957020112577174
$data[9,3] = 0.9867934623294556
$data[9,4] = 0.9250863404092976
$data[9,5] = 1.092331723448326
$data[9,6] = 0.8684472120320456
$data[9,7] = 0.9544445180985303
$data[9,8] = 0.957020112577174
$data[9,9] = 0.9719067874533148
$data[9,10] = 0.9484965639313062
$data[9,11] = 0.9640205614824716
$ws.Range("B10:M19").Value2 = $data

# A column (row index) for the three brand-new rows 17-19
$ws.Range("A17").Value2 = 15
$ws.Range("A18").Value2 = 16
$ws.Range("A19").Value2 = 17

# Copy the existing bold/centered/bordered style from A16 onto the new A17:A19 cells
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
